$wb = $excel.ActiveWorkbook

# Update the existing "prompt" sheet's text value (shared string content stays the
# same logical text - "This is updated english prompt" - it just ends up being
# re-indexed in sharedStrings.xml once the old, now-unused "This is english prompt"
# string is dropped).
$promptSheet = $wb.Worksheets.Item("prompt")
$promptSheet.Range("A1").Value2 = "This is updated english prompt"

# Add a new worksheet named "survey" after the last existing sheet ("prompt").
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$surveySheet = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$surveySheet.Name = "survey"
$surveySheet.Range("A1").Value2 = "A Updated Health Survey"

# Make the new "survey" sheet the active / selected tab.
$surveySheet.Select()
